$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 36, pushing existing rows 36-130 down to 37-131.
$ws.Rows.Item(36).Insert(-4121)

# Populate the newly inserted row 36 with the new weekly data point.
$ws.Range("A36").Value = 9
$ws.Range("B36").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C36").Value = "Metropolitana"
$ws.Range("D36").Value = 45133
$ws.Range("E36").Value = 13
$ws.Range("F36").Value = 100112005
$ws.Range("G36").Value = "Puerro"
$ws.Range("H36").Value = "Sin especificar"
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 70
$ws.Range("K36").Value = 8000
$ws.Range("L36").Value = 8000
$ws.Range("M36").Value = 8000
$ws.Range("N36").Value = '$/paquete 20 unidades'
$ws.Range("O36").Value = "Provincia de Chacabuco"
$ws.Range("P36").Value = 400
$ws.Range("Q36").Value = 20
$ws.Range("R36").Value = "Hortaliza"
